$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - same value/formatting as the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data column H2:H5 ("Save" flag values)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
